# Generate Report for Handoff
#
# Updates the localization-status report so the "Handed back" row now
# reflects that a fresh handoff package is ready:
#   - Status text changes from "Handed back: in sync with en-US"
#     to "Ready for handoff" (Overview!E2/F2, zh-cn!C2, de-de!C2)
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are bumped to the new handoff time
#   - The Status column narrows to fit the shorter text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-09-07 01:19:00"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-09-07 01:18:56"

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-09-07 01:19:00"

# --- Shrink the now-narrower Status columns to fit the new text --------
# (Overview columns E & F hold the per-language status; zh-cn/de-de
# column C is the Status column.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
